# Correção nos dados e início da análise PNAD 2009
#
# The sheet had extra "category header" rows (sexo, cor ou raça, grupos de
# idade, nível de instrução, classes de rendimento mensal domiciliar per
# capita) that carried no numeric data, plus two trailing footnote rows
# (fonte / notes) at the very end. Remove them entirely; Excel's row
# deletion naturally shifts the remaining data rows up so each label keeps
# its correct coefficient-of-variation values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so row numbers of not-yet-deleted rows don't shift.
$ws.Rows.Item(35).Delete()
$ws.Rows.Item(34).Delete()
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
